$d = $word.ActiveDocument

$replacements = @(
    @("820÷3=273, 1", "430÷8=53, 6"),
    @("713÷5=142, 3", "152÷3=50, 2"),
    @("586÷7=83, 5", "314÷6=52, 2"),
    @("826÷9=91, 7", "715÷2=357, 1"),
    @("526÷2=263, 0", "442÷6=73, 4"),
    @("927÷4=231, 3", "378÷7=54, 0"),
    @("444÷8=55, 4", "931÷9=103, 4"),
    @("129÷4=32, 1", "981÷6=163, 3"),
    @("228÷7=32, 4", "684÷8=85, 4"),
    @("806÷2=403, 0", "464÷7=66, 2"),
    @("698÷4=174, 2", "573÷7=81, 6"),
    @("888÷4=222, 0", "554÷4=138, 2"),
    @("688÷6=114, 4", "497÷9=55, 2"),
    @("126÷6=21, 0", "930÷2=465, 0"),
    @("454÷7=64, 6", "643÷3=214, 1"),
    @("663÷5=132, 3", "745÷6=124, 1"),
    @("845÷2=422, 1", "495÷3=165, 0"),
    @("127÷8=15, 7", "404÷7=57, 5"),
    @("249÷5=49, 4", "900÷6=150, 0"),
    @("999÷6=166, 3", "267÷6=44, 3"),
    @("228÷6=38, 0", "581÷9=64, 5"),
    @("836÷5=167, 1", "489÷9=54, 3"),
    @("317÷4=79, 1", "691÷6=115, 1"),
    @("270÷5=54, 0", "547÷6=91, 1"),
    @("185÷2=92, 1", "702÷3=234, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Host "Done applying replacements"
